# Adds 5 new words to the "Przerobione" sheet (rows 62-66) and refreshes
# the "5 losowych" sheet so it mirrors those newly added 5 words.

$wb = $excel.ActiveWorkbook

$przerobione = $wb.Worksheets.Item("Przerobione")
$losowych = $wb.Worksheets.Item("5 losowych")

# New vocabulary entries: id, hanzi, pinyin, english
$newWords = @(
    @(710, "然而", "rán'ér", "however"),
    @(292, "工资", "gōngzī", "wages"),
    @(848, "疼", "téng", "hurt"),
    @(1101, "责任", "zérèn", "responsibility"),
    @(597, "米", "mǐ", "rice")
)

$startRow = 62
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $word = $newWords[$i]

    $przerobione.Cells.Item($row, 1).Value = $word[0]
    $przerobione.Cells.Item($row, 2).Value = $word[1]
    $przerobione.Cells.Item($row, 3).Value = $word[2]
    $przerobione.Cells.Item($row, 4).Value = $word[3]
}

# "5 losowych" (5 random) mirrors the latest 5 words added to "Przerobione"
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = 2 + $i
    $word = $newWords[$i]

    $losowych.Cells.Item($row, 1).Value = $word[0]
    $losowych.Cells.Item($row, 2).Value = $word[1]
    $losowych.Cells.Item($row, 3).Value = $word[2]
    $losowych.Cells.Item($row, 4).Value = $word[3]
}
